# Updated cryptos list on Mon Sep 25 07:14:09 UTC 2023 with GitHub Actions
#
# For each data row, refresh the Price (column D) and Volume(1h) (column E)
# cells. A handful of Price values are plain decimal-looking numbers
# ("209.68", "0.498", ...); Excel's COM type-coercion would otherwise store
# those as numeric cells, but the source data keeps them as text (same as
# every other row), so for those we flip the cell to Text format first -
# exactly what typing `'209.68` into Excel does - before writing the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# rows whose new Price value reads as a plain number and must be pinned to Text
$priceNeedsTextFormat = @(5,6,10,15,16,19,21,23,24,25,27,29,30,32,37,41,42,44,46,49,51)

# row -> (new Price, new Volume(1h)) ; $null means "leave this column alone"
$updates = @{
    2  = @("26.190.24", "  -1.98%  ")
    3  = @("1.582.36",  "  -1.21%  ")
    4  = @($null,       "  -0.34%  ")
    5  = @("209.68",    "  -0.99%  ")
    6  = @("0.498",     "  -2.87%  ")
    7  = @($null,       "  -0.32%  ")
    8  = @($null,       "  -1.51%  ")
    9  = @($null,       "  -0.72%  ")
    10 = @("19.49",     "  -1.32%  ")
    11 = @($null,       "  +0.32%  ")
    12 = @("1.805.11",  "  -1.18%  ")
    13 = @("1.605.15",  "  +0.41%  ")
    14 = @($null,       "  -0.14%  ")
    15 = @("0.516",     "  -1.45%  ")
    16 = @("64.48",     "  -0.74%  ")
    17 = @("26.197.64", "  -1.84%  ")
    18 = @("0.0₃0733",  "  -1.18%  ")
    19 = @("7.25",      "  +1.03%  ")
    20 = @($null,       "  -0.31%  ")
    21 = @("206.39",    "  -1.99%  ")
    22 = @($null,       "  -0.67%  ")
    23 = @("2.20",      "  -3.44%  ")
    24 = @("8.86",      "  -1.08%  ")
    25 = @("144.87",    "  +0.67%  ")
    26 = @($null,       "  -0.34%  ")
    27 = @("7.02",      "  -0.88%  ")
    28 = @($null,       "  -1.15%  ")
    29 = @("15.21",     "  -1.02%  ")
    30 = @("0.0504",    "  -1.43%  ")
    31 = @($null,       "  -0.75%  ")
    32 = @("3.22",      "  -1.46%  ")
    33 = @($null,       "  -1.07%  ")
    34 = @("1.282.16",  "  -0.82%  ")
    35 = @($null,       "  +7.89%  ")
    37 = @("0.604",     "  +0.93%  ")
    38 = @($null,       "  -1.06%  ")
    39 = @($null,       "  -1.48%  ")
    40 = @($null,       "  -1.83%  ")
    41 = @("5.56",      "  +3.12%  ")
    42 = @("0.768",     "  -1.61%  ")
    43 = @($null,       "  -2.80%  ")
    44 = @("62.21",     "  -1.39%  ")
    45 = @("1.718.41",  $null)
    46 = @("88.76",     "  -2.06%  ")
    47 = @($null,       "  -0.46%  ")
    48 = @($null,       "  -0.96%  ")
    49 = @("0.0506",    "  -1.94%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]

    if ($null -ne $newPrice) {
        if ($priceNeedsTextFormat -contains $row) {
            Set-TextValue $row 4 $newPrice
        } else {
            $ws.Cells.Item($row, 4).Value = $newPrice
        }
    }

    if ($null -ne $newVolume) {
        $ws.Cells.Item($row, 5).Value = $newVolume
    }
}

# Row 50 used to be USDD and row 51 used to be EnergySwap. EnergySwap dropped
# out of the list; USDD shifted down a slot and BabyDogeCoin entered at 50.
$ws.Cells.Item(50, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 50 4 "0.0₇0967"
$ws.Cells.Item(50, 5).Value = "  -9.10%  "

$ws.Cells.Item(51, 2).Value = "USDD"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue 51 4 "1.00"
$ws.Cells.Item(51, 5).Value = "  -0.19%  "
